$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assigns a text value to a cell, preserving it as text even when the
# string looks like a number (Excel would otherwise auto-convert it to a float).
function Set-TextValue($sheet, $cellRef, $text) {
    $cell = $sheet.Range($cellRef)
    if ($text -match '^[+-]?[0-9]*\.?[0-9]+$') {
        $cell.Value = "'" + $text
    } else {
        $cell.Value = $text
    }
}

# Row 2
Set-TextValue $ws "D2" "97.291.51"
Set-TextValue $ws "E2" "  +1.71%  "

# Row 3
Set-TextValue $ws "D3" "3.595.09"
Set-TextValue $ws "E3" "  -0.94%  "

# Row 4
Set-TextValue $ws "E4" "  -0.03%  "

# Row 5
Set-TextValue $ws "D5" "244.78"
Set-TextValue $ws "E5" "  +3.21%  "

# Row 6
Set-TextValue $ws "D6" "655.29"
Set-TextValue $ws "E6" "  -0.59%  "

# Row 7
Set-TextValue $ws "E7" "  +12.86%  "

# Row 8
Set-TextValue $ws "D8" "0.415"
Set-TextValue $ws "E8" "  +3.08%  "

# Row 9
Set-TextValue $ws "D9" "1.07"
Set-TextValue $ws "E9" "  +6.67%  "

# Row 10
Set-TextValue $ws "D10" "0.999"
Set-TextValue $ws "E10" "  -0.03%  "

# Row 11
Set-TextValue $ws "D11" "3.595.28"
Set-TextValue $ws "E11" "  -0.90%  "

# Row 12
Set-TextValue $ws "D12" "43.89"
Set-TextValue $ws "E12" "  +3.70%  "

# Row 13
Set-TextValue $ws "E13" "  +1.53%  "

# Row 14
Set-TextValue $ws "D14" "6.44"

# Row 15
Set-TextValue $ws "D15" "4.261.42"
Set-TextValue $ws "E15" "  -1.01%  "

# Row 16
Set-TextValue $ws "D16" "96.933.93"
Set-TextValue $ws "E16" "  +1.52%  "

# Row 17
Set-TextValue $ws "D17" "0.0000258"
Set-TextValue $ws "E17" "  +1.80%  "

# Row 18
Set-TextValue $ws "D18" "3.595.27"
Set-TextValue $ws "E18" "  -0.82%  "

# Row 19
Set-TextValue $ws "B19" "Uniswap"
Set-TextValue $ws "C19" "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue $ws "D19" "12.76"
Set-TextValue $ws "E19" "  -1.47%  "

# Row 20
Set-TextValue $ws "B20" "Polkadot"
Set-TextValue $ws "C20" "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue $ws "D20" "7.76"
Set-TextValue $ws "E20" "  -2.22%  "

# Row 21
Set-TextValue $ws "D21" "18.14"
Set-TextValue $ws "E21" "  +0.59%  "

# Row 22
Set-TextValue $ws "E22" "  +9.70%  "

# Row 23
Set-TextValue $ws "B23" "BitcoinCash"
Set-TextValue $ws "C23" "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue $ws "D23" "510.91"
Set-TextValue $ws "E23" "  +1.15%  "

# Row 24
Set-TextValue $ws "B24" "SuiNetwork"
Set-TextValue $ws "C24" "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
Set-TextValue $ws "D24" "3.43"
Set-TextValue $ws "E24" "  -3.39%  "

# Row 25
Set-TextValue $ws "E25" "  +2.43%  "

# Row 26
Set-TextValue $ws "D26" "6.94"
Set-TextValue $ws "E26" "  +4.26%  "

# Row 27
Set-TextValue $ws "D27" "97.27"
Set-TextValue $ws "E27" "  +5.79%  "

# Row 28
Set-TextValue $ws "D28" "13.22"
Set-TextValue $ws "E28" "  +5.42%  "

# Row 29
Set-TextValue $ws "D29" "3.787.12"
Set-TextValue $ws "E29" "  -0.98%  "

# Row 30
Set-TextValue $ws "E30" "  -1.68%  "

# Row 31
Set-TextValue $ws "E31" "  +9.49%  "

# Row 32
Set-TextValue $ws "D32" "11.59"
Set-TextValue $ws "E32" "  +2.60%  "

# Row 33
Set-TextValue $ws "D33" "1.00"
Set-TextValue $ws "E33" "  -0.02%  "

# Row 34
Set-TextValue $ws "D34" "0.186"
Set-TextValue $ws "E34" "  +5.16%  "

# Row 35
Set-TextValue $ws "E35" "  +0.46%  "

# Row 36
Set-TextValue $ws "D36" "31.63"
Set-TextValue $ws "E36" "  -2.40%  "

# Row 37
Set-TextValue $ws "D37" "8.87"
Set-TextValue $ws "E37" "  +9.14%  "

# Row 38
Set-TextValue $ws "D38" "623.06"
Set-TextValue $ws "E38" "  +9.49%  "

# Row 39
Set-TextValue $ws "D39" "0.573"
Set-TextValue $ws "E39" "  +2.09%  "

# Row 40
Set-TextValue $ws "E40" "  +11.62%  "

# Row 41
Set-TextValue $ws "D41" "0.152"
Set-TextValue $ws "E41" "  +1.31%  "

# Row 42
Set-TextValue $ws "E42" "  +0.02%  "

# Row 43
Set-TextValue $ws "D43" "0.918"
Set-TextValue $ws "E43" "  +0.08%  "

# Row 44
Set-TextValue $ws "D44" "1.86"
Set-TextValue $ws "E44" "  +6.31%  "

# Row 45
Set-TextValue $ws "B45" "VeChain"
Set-TextValue $ws "C45" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws "D45" "0.0434"
Set-TextValue $ws "E45" "  +5.00%  "

# Row 46
Set-TextValue $ws "B46" "Filecoin"
Set-TextValue $ws "C46" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws "D46" "5.83"
Set-TextValue $ws "E46" "  +3.04%  "

# Row 47
Set-TextValue $ws "E47" "  +1.71%  "

# Row 48
Set-TextValue $ws "D48" "23.57"
Set-TextValue $ws "E48" "  -0.41%  "

# Row 49
Set-TextValue $ws "D49" "33.44"
Set-TextValue $ws "E49" "  -8.78%  "

# Row 50
Set-TextValue $ws "E50" "  -0.07%  "

# Row 51
Set-TextValue $ws "D51" "8.34"
Set-TextValue $ws "E51" "  +3.87%  "
